$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1510.1333
$ws.Range("J41").Value = 2148.4
$ws.Range("L41").Value = 2148.4
$ws.Range("N41").Value = -3028.4

$ws.Range("H53").Value = 1461.5333
$ws.Range("I53").Value = 1043.4546
$ws.Range("K53").Value = 1043.4546
$ws.Range("M53").Value = -406.4546

$ws.Range("H69").Value = 3752
$ws.Range("I69").Value = 3561.6667
$ws.Range("J69").Value = 3894.75
$ws.Range("K69").Value = 10685.0001
$ws.Range("L69").Value = 11684.25
$ws.Range("M69").Value = -9811.000100000001
$ws.Range("N69").Value = -13432.25

$ws.Range("H72").Value = 3752
$ws.Range("I72").Value = 3561.6667
$ws.Range("J72").Value = 3894.75
$ws.Range("K72").Value = 32055.0003
$ws.Range("L72").Value = 35052.75
$ws.Range("M72").Value = -27687.0003
$ws.Range("N72").Value = -43788.75

$ws.Range("H76").Value = 3578.1667
$ws.Range("I76").Value = 3567
$ws.Range("J76").Value = 3589.3333
$ws.Range("K76").Value = 3567
$ws.Range("L76").Value = 3589.3333
$ws.Range("M76").Value = -3252
$ws.Range("N76").Value = -4219.3333

$ws.Range("H79").Value = 3578.1667
$ws.Range("I79").Value = 3567
$ws.Range("J79").Value = 3589.3333
$ws.Range("K79").Value = 3567
$ws.Range("L79").Value = 3589.3333
$ws.Range("M79").Value = -2475
$ws.Range("N79").Value = -5773.3333

$ws.Range("H137").Value = 1177.4736
$ws.Range("I137").Value = 1154.125
$ws.Range("J137").Value = 1302
$ws.Range("K137").Value = 3462.375
$ws.Range("L137").Value = 3906
$ws.Range("M137").Value = -912.375
$ws.Range("N137").Value = -9006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2018.1111
$ws.Range("I61").Value = 1166.2858
$ws.Range("J61").Value = 4999.5
$ws.Range("K61").Value = 1166.2858
$ws.Range("L61").Value = 4999.5
$ws.Range("M61").Value = -954.2858000000001
$ws.Range("N61").Value = -5423.5

$ws.Range("H97").Value = 354.85715
$ws.Range("I97").Value = 300.9
$ws.Range("J97").Value = 489.75
$ws.Range("K97").Value = 300.9
$ws.Range("L97").Value = 489.75
$ws.Range("M97").Value = 195.1
$ws.Range("N97").Value = -1481.75

$ws.Range("H122").Value = 1425.8334
$ws.Range("I122").Value = 1425.8334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4277.5002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1827.5002
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3125.4285
$ws.Range("I132").Value = 3067.5454
$ws.Range("J132").Value = 3337.6667
$ws.Range("K132").Value = 9202.636200000001
$ws.Range("L132").Value = 10013.0001
$ws.Range("M132").Value = -6672.636200000001
$ws.Range("N132").Value = -15073.0001

$ws.Range("H136").Value = 2018.1111
$ws.Range("I136").Value = 1166.2858
$ws.Range("J136").Value = 4999.5
$ws.Range("K136").Value = 3498.8574
$ws.Range("L136").Value = 14998.5
$ws.Range("M136").Value = -948.8574000000003
$ws.Range("N136").Value = -20098.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1675.0869
$ws.Range("I20").Value = 1650.4706
$ws.Range("J20").Value = 1744.8334
$ws.Range("K20").Value = 1650.4706
$ws.Range("L20").Value = 1744.8334
$ws.Range("M20").Value = -1403.4706
$ws.Range("N20").Value = -2238.8334

$ws.Range("H105").Value = 111113704
$ws.Range("I105").Value = 125002540
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 125002540
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -125000793
$ws.Range("N105").Value = -6505

$ws.Range("H134").Value = 8104.375
$ws.Range("I134").Value = 1722.5834
$ws.Range("J134").Value = 27249.75
$ws.Range("K134").Value = 5167.7502
$ws.Range("L134").Value = 81749.25
$ws.Range("M134").Value = -2632.7502
$ws.Range("N134").Value = -86819.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 111112424
$ws.Range("I16").Value = 166667890
$ws.Range("J16").Value = 1492.6666
$ws.Range("K16").Value = 166667890
$ws.Range("L16").Value = 1492.6666
$ws.Range("M16").Value = -166667603
$ws.Range("N16").Value = -2066.6666

$ws.Range("H31").Value = 2191.862
$ws.Range("I31").Value = 1072.4166
$ws.Range("J31").Value = 2982.0588
$ws.Range("K31").Value = 1072.4166
$ws.Range("L31").Value = 2982.0588
$ws.Range("M31").Value = -777.4166
$ws.Range("N31").Value = -3572.0588

$ws.Range("H34").Value = 2191.862
$ws.Range("I34").Value = 1072.4166
$ws.Range("J34").Value = 2982.0588
$ws.Range("K34").Value = 1072.4166
$ws.Range("L34").Value = 2982.0588
$ws.Range("M34").Value = -870.4166
$ws.Range("N34").Value = -3386.0588

$ws.Range("H69").Value = 5000
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251

$ws.Range("H72").Value = 5000
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256

$ws.Range("H107").Value = 623.4783
$ws.Range("I107").Value = 429.41666
$ws.Range("K107").Value = 429.41666
$ws.Range("M107").Value = 1490.58334

$ws.Range("H113").Value = 111112424
$ws.Range("I113").Value = 166667890
$ws.Range("J113").Value = 1492.6666
$ws.Range("K113").Value = 166667890
$ws.Range("L113").Value = 1492.6666
$ws.Range("M113").Value = -166665720
$ws.Range("N113").Value = -5832.6666

$ws.Range("H141").Value = 32732.727
$ws.Range("J141").Value = 32732.727
$ws.Range("L141").Value = 32732.727
$ws.Range("N141").Value = -43092.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 63.75
$ws.Range("I2").Value = 47.5
$ws.Range("K2").Value = 285
$ws.Range("M2").Value = -172

$ws.Range("H61").Value = 201
$ws.Range("I61").Value = 121.4
$ws.Range("J61").Value = 400
$ws.Range("K61").Value = 364.2
$ws.Range("L61").Value = 1200
$ws.Range("M61").Value = -149.2
$ws.Range("N61").Value = -1630

$ws.Range("H69").Value = 2158.6316
$ws.Range("I69").Value = 600
$ws.Range("K69").Value = 1800
$ws.Range("M69").Value = -989

$ws.Range("H72").Value = 2158.6316
$ws.Range("I72").Value = 600
$ws.Range("K72").Value = 5400
$ws.Range("M72").Value = -1344

$ws.Range("H139").Value = 3754.875
$ws.Range("I139").Value = 3754.875
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 11264.625
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -6124.625
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 1869.4783
$ws.Range("I140").Value = 1822.6364
$ws.Range("K140").Value = 5467.9092
$ws.Range("M140").Value = -287.9092000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 5000
$ws.Range("J47").Value = 5000
$ws.Range("L47").Value = 5000
$ws.Range("N47").Value = -6136

$ws.Range("H122").Value = 251617.83
$ws.Range("I122").Value = 1069
$ws.Range("J122").Value = 502166.66
$ws.Range("K122").Value = 3207
$ws.Range("L122").Value = 1506499.98
$ws.Range("M122").Value = -757
$ws.Range("N122").Value = -1511399.98

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 13891456
$ws.Range("I122").Value = 20835708
$ws.Range("J122").Value = 2952.3333
$ws.Range("K122").Value = 62507124
$ws.Range("L122").Value = 8856.999899999999
$ws.Range("M122").Value = -62504674
$ws.Range("N122").Value = -13756.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 29418636
$ws.Range("I62").Value = 33338720
$ws.Range("J62").Value = 18000
$ws.Range("K62").Value = 33338720
$ws.Range("L62").Value = 18000
$ws.Range("M62").Value = -33338096
$ws.Range("N62").Value = -19248

$ws.Range("H65").Value = 29418636
$ws.Range("I65").Value = 33338720
$ws.Range("J65").Value = 18000
$ws.Range("K65").Value = 166693600
$ws.Range("L65").Value = 90000
$ws.Range("M65").Value = -166690480
$ws.Range("N65").Value = -96240

$ws.Range("H122").Value = 42001904
$ws.Range("I122").Value = 48463600
$ws.Range("J122").Value = 875
$ws.Range("K122").Value = 145390800
$ws.Range("L122").Value = 2625
$ws.Range("M122").Value = -145388350
$ws.Range("N122").Value = -7525

$ws.Range("H132").Value = 4900.7334
$ws.Range("I132").Value = 4845.8887
$ws.Range("K132").Value = 14537.6661
$ws.Range("M132").Value = -12007.6661

$ws.Range("H136").Value = 975.75
$ws.Range("I136").Value = 975.75
$ws.Range("K136").Value = 2927.25
$ws.Range("M136").Value = -377.25
